$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 223.6923
$ws.Cells.Item(5, 9).Value = 177.25
$ws.Cells.Item(5, 10).Value = 298
$ws.Cells.Item(5, 11).Value = 177.25
$ws.Cells.Item(5, 12).Value = 298
$ws.Cells.Item(5, 13).Value = -62.25
$ws.Cells.Item(5, 14).Value = -528

$ws.Cells.Item(40, 8).Value = 1880.8966
$ws.Cells.Item(40, 9).Value = 1643.75
$ws.Cells.Item(40, 10).Value = 2172.7693
$ws.Cells.Item(40, 11).Value = 1643.75
$ws.Cells.Item(40, 12).Value = 2172.7693
$ws.Cells.Item(40, 13).Value = -1468.75
$ws.Cells.Item(40, 14).Value = -2522.7693

$ws.Cells.Item(41, 8).Value = 284.6
$ws.Cells.Item(41, 9).Value = 286
$ws.Cells.Item(41, 10).Value = 282.5
$ws.Cells.Item(41, 11).Value = 286
$ws.Cells.Item(41, 12).Value = 282.5
$ws.Cells.Item(41, 13).Value = 154
$ws.Cells.Item(41, 14).Value = -1162.5

$ws.Cells.Item(45, 8).Value = 3008.5
$ws.Cells.Item(45, 9).Value = 1017
$ws.Cells.Item(45, 11).Value = 3051
$ws.Cells.Item(45, 13).Value = -2859

$ws.Cells.Item(54, 8).Value = 10724.75
$ws.Cells.Item(54, 9).Value = 8166.3335
$ws.Cells.Item(54, 10).Value = 18400
$ws.Cells.Item(54, 11).Value = 8166.3335
$ws.Cells.Item(54, 12).Value = 18400
$ws.Cells.Item(54, 13).Value = -7680.3335
$ws.Cells.Item(54, 14).Value = -19372

$ws.Cells.Item(61, 8).Value = 120.666664
$ws.Cells.Item(61, 9).Value = 120.666664
$ws.Cells.Item(61, 11).Value = 361.999992
$ws.Cells.Item(61, 13).Value = -189.999992

$ws.Cells.Item(116, 8).Value = 7185.2915
$ws.Cells.Item(116, 9).Value = 5323.7393
$ws.Cells.Item(116, 11).Value = 5323.7393
$ws.Cells.Item(116, 13).Value = -1881.7393

$ws.Cells.Item(132, 8).Value = 3751.6943
$ws.Cells.Item(132, 9).Value = 1507.2916
$ws.Cells.Item(132, 10).Value = 8240.5
$ws.Cells.Item(132, 11).Value = 4521.8748
$ws.Cells.Item(132, 12).Value = 24721.5
$ws.Cells.Item(132, 13).Value = -1991.8748
$ws.Cells.Item(132, 14).Value = -29781.5

$ws.Cells.Item(135, 8).Value = 928.6842
$ws.Cells.Item(135, 9).Value = 923.9286
$ws.Cells.Item(135, 10).Value = 942
$ws.Cells.Item(135, 11).Value = 8315.357399999999
$ws.Cells.Item(135, 12).Value = 8478
$ws.Cells.Item(135, 13).Value = -5780.357399999999
$ws.Cells.Item(135, 14).Value = -13548

$ws.Cells.Item(139, 8).Value = 36666.332
$ws.Cells.Item(139, 10).Value = 45000
$ws.Cells.Item(139, 12).Value = 45000
$ws.Cells.Item(139, 14).Value = -55280

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 3608
$ws.Cells.Item(3, 9).Value = 1416.6666
$ws.Cells.Item(3, 11).Value = 1416.6666
$ws.Cells.Item(3, 13).Value = -1301.6666

$ws.Cells.Item(94, 8).Value = 29147.5
$ws.Cells.Item(94, 10).Value = 29147.5
$ws.Cells.Item(94, 12).Value = 29147.5
$ws.Cells.Item(94, 14).Value = -30949.5

$ws.Cells.Item(110, 8).Value = 1810.5454
$ws.Cells.Item(110, 9).Value = 1921.6
$ws.Cells.Item(110, 11).Value = 1921.6
$ws.Cells.Item(110, 13).Value = 123.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 5289.125
$ws.Cells.Item(75, 9).Value = 5289.125
$ws.Cells.Item(75, 11).Value = 5289.125
$ws.Cells.Item(75, 13).Value = -4353.125

$ws.Cells.Item(78, 8).Value = 5289.125
$ws.Cells.Item(78, 9).Value = 5289.125
$ws.Cells.Item(78, 11).Value = 15867.375
$ws.Cells.Item(78, 13).Value = -11187.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 31999.6
$ws.Cells.Item(20, 10).Value = 31999.6
$ws.Cells.Item(20, 12).Value = 31999.6
$ws.Cells.Item(20, 14).Value = -32471.6

$ws.Cells.Item(30, 8).Value = 31999.6
$ws.Cells.Item(30, 10).Value = 31999.6
$ws.Cells.Item(30, 12).Value = 31999.6
$ws.Cells.Item(30, 14).Value = -32181.6

$ws.Cells.Item(127, 8).Value = 31804.625
$ws.Cells.Item(127, 10).Value = 31804.625
$ws.Cells.Item(127, 12).Value = 31804.625
$ws.Cells.Item(127, 14).Value = -41724.625

$ws.Cells.Item(128, 8).Value = 31999.6
$ws.Cells.Item(128, 10).Value = 31999.6
$ws.Cells.Item(128, 12).Value = 31999.6
$ws.Cells.Item(128, 14).Value = -41959.6

$ws.Cells.Item(131, 8).Value = 29307.889
$ws.Cells.Item(131, 10).Value = 29307.889
$ws.Cells.Item(131, 12).Value = 29307.889
$ws.Cells.Item(131, 14).Value = -39387.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 3538.125
$ws.Cells.Item(51, 10).Value = 3640.6667
$ws.Cells.Item(51, 12).Value = 10922.0001
$ws.Cells.Item(51, 14).Value = -11842.0001

$ws.Cells.Item(131, 8).Value = 891.3838500000001
$ws.Cells.Item(131, 9).Value = 609.8570999999999
$ws.Cells.Item(131, 10).Value = 912.8043
$ws.Cells.Item(131, 11).Value = 1829.5713
$ws.Cells.Item(131, 12).Value = 2738.4129
$ws.Cells.Item(131, 13).Value = 3210.4287
$ws.Cells.Item(131, 14).Value = -12818.4129

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 984.5
$ws.Cells.Item(9, 9).Value = 984.5
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 984.5
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = -814.5
$ws.Cells.Item(9, 14).ClearContents()

$ws.Cells.Item(97, 8).Value = 1877.75
$ws.Cells.Item(97, 9).Value = 1700
$ws.Cells.Item(97, 10).Value = 1937
$ws.Cells.Item(97, 11).Value = 1700
$ws.Cells.Item(97, 12).Value = 1937
$ws.Cells.Item(97, 13).Value = -1204
$ws.Cells.Item(97, 14).Value = -2929

$ws.Cells.Item(113, 8).Value = 1351.6666
$ws.Cells.Item(113, 9).Value = 833.5
$ws.Cells.Item(113, 10).Value = 1869.8334
$ws.Cells.Item(113, 11).Value = 833.5
$ws.Cells.Item(113, 12).Value = 1869.8334
$ws.Cells.Item(113, 13).Value = 1336.5
$ws.Cells.Item(113, 14).Value = -6209.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(57, 8).Value = 20015.334
$ws.Cells.Item(57, 9).Value = 20000
$ws.Cells.Item(57, 10).Value = 20023
$ws.Cells.Item(57, 11).Value = 20000
$ws.Cells.Item(57, 12).Value = 20023
$ws.Cells.Item(57, 13).Value = -19434
$ws.Cells.Item(57, 14).Value = -21155

$ws.Cells.Item(93, 8).Value = 175675.75
$ws.Cells.Item(93, 9).Value = 850
$ws.Cells.Item(93, 10).Value = 350501.5
$ws.Cells.Item(93, 11).Value = 850
$ws.Cells.Item(93, 12).Value = 350501.5
$ws.Cells.Item(93, 13).Value = 398
$ws.Cells.Item(93, 14).Value = -352997.5

$ws.Cells.Item(100, 8).Value = 2571.182
$ws.Cells.Item(100, 9).Value = 2328.3
$ws.Cells.Item(100, 11).Value = 2328.3
$ws.Cells.Item(100, 13).Value = -1787.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1108.1818
$ws.Cells.Item(81, 9).Value = 1058
$ws.Cells.Item(81, 10).Value = 1150
$ws.Cells.Item(81, 11).Value = 2116
$ws.Cells.Item(81, 12).Value = 2300
$ws.Cells.Item(81, 13).Value = -1055
$ws.Cells.Item(81, 14).Value = -4422

$ws.Cells.Item(84, 8).Value = 1108.1818
$ws.Cells.Item(84, 9).Value = 1058
$ws.Cells.Item(84, 10).Value = 1150
$ws.Cells.Item(84, 11).Value = 10580
$ws.Cells.Item(84, 12).Value = 11500
$ws.Cells.Item(84, 13).Value = -5276
$ws.Cells.Item(84, 14).Value = -22108

$ws.Cells.Item(126, 8).Value = 52632290
$ws.Cells.Item(126, 9).Value = 83333970
$ws.Cells.Item(126, 10).Value = 839.5714
$ws.Cells.Item(126, 11).Value = 250001910
$ws.Cells.Item(126, 12).Value = 2518.7142
$ws.Cells.Item(126, 13).Value = -249999440
$ws.Cells.Item(126, 14).Value = -7458.7142

$ws.Cells.Item(132, 8).Value = 3335644
$ws.Cells.Item(132, 9).Value = 4547769
$ws.Cells.Item(132, 11).Value = 13643307
$ws.Cells.Item(132, 13).Value = -13640777
